$wb = $excel.ActiveWorkbook

# xlEdgeLeft=7, xlEdgeTop=8, xlEdgeBottom=9, xlEdgeRight=10, xlContinuous=1
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlContinuous = 1

function Set-InnerHeaderBorder($range) {
    # Style used by the "inner" cell(s) of a merged header: only a
    # top+bottom box (matches the pre-existing borderId 4 in styles.xml).
    $range.ClearFormats()
    $range.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
    $range.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
}

function Set-RightHeaderBorder($range) {
    # Style used by the rightmost cell of a merged header: top+bottom+right
    # box (matches the pre-existing borderId 5 in styles.xml). Setting the
    # right edge first keeps this from ever passing through the exact
    # "top+bottom only" combination used by Set-InnerHeaderBorder above,
    # which would otherwise get it (temporarily) assigned to the very same
    # style record and leave a stray/duplicate entry behind in cellXfs.
    $range.ClearFormats()
    $range.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
    $range.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
    $range.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
}

# --- Sheet 1: quality_comparison ---------------------------------------
$ws1 = $wb.Worksheets.Item("quality_comparison")

Set-InnerHeaderBorder $ws1.Range("C1")
Set-RightHeaderBorder $ws1.Range("D1")

$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ----------------------------------
$ws2 = $wb.Worksheets.Item("computational_comparison")

Set-InnerHeaderBorder $ws2.Range("C1")
Set-RightHeaderBorder $ws2.Range("D1")
Set-InnerHeaderBorder $ws2.Range("F1")
Set-RightHeaderBorder $ws2.Range("G1")

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# G5 becomes a genuinely empty cell (the element disappears entirely).
$ws2.Range("G5").ClearContents()

Write-Output "edits applied"
